# xpath logic and code restructure
# Insert two new articles into the NY Times articles sheet:
#   - "Golf's Big Deal Veers Off Course" goes in ahead of the existing
#     Taylor Swift row (new row 3, pushing Taylor Swift down to row 4).
#   - "Richard Lyons..." and "36 Hours in Toronto" are inserted ahead of
#     the J. Cole row (new rows 5-6, pushing J. Cole and the Beyonce
#     Popcast rows down to rows 7-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a single row above current row 3 (Taylor Swift / April 11),
# making room for the new Golf article.
$ws.Rows.Item(3).Insert()

# Insert two rows above current row 5 (J. Cole / April 8), making room
# for the Richard Lyons and Toronto articles.
$ws.Rows.Item(5).Resize(2).Insert()

# New row 3: Golf's Big Deal Veers Off Course
$ws.Range("A3").Value = "Golf’s Big Deal Veers Off Course"
$ws.Range("B3").Value = "The Masters tournament should be all about sport, but the unresolved fight between the PGA Tour and LIV Golf looms over the competition."
$ws.Range("C3").Value = "April 13"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = "images/DB13-golf-jhfq-threeByTwoSmallAt2X.jpg"

# New row 5: Richard Lyons, Former Business School Dean, Will Be U.C. Berkeley's New Chancellor
$ws.Range("A5").Value = "Richard Lyons, Former Business School Dean, Will Be U.C. Berkeley’s New Chancellor"
$ws.Range("B5").Value = "The appointment comes as Berkeley and college campuses across the country are facing turmoil over free speech, racial and political diversity, and affordability."
$ws.Range("C5").Value = "April 11"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = "images/11californiatoday-ucb-threeByTwoSmallAt2X.jpg"

# New row 6: 36 Hours in Toronto
$ws.Range("A6").Value = "36 Hours in Toronto"
$ws.Range("B6").Value = "Savor the diversity of this lakefront city though its hidden bars, small-but-fascinating museums and restaurants with dishes like jerk chicken chow mein and Hong Kong-style French toast."
$ws.Range("C6").Value = "April 11"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = "images/things-to-do-toronto-01-hwlf-threeByTwoSmallAt2X.jpg"
